$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045522924504251
$ws.Range("D2").Value = 1.046990129193348
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.044360905294209
$ws.Range("I2").Value = 1.041598837292934
$ws.Range("J2").Value = 1.050582205824155
$ws.Range("K2").Value = 1.049754019732186
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.047132172547802
$ws.Range("N2").Value = 1.052074153059889

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046885353185557
$ws.Range("D3").Value = 1.048052940320931
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.046349716357448
$ws.Range("I3").Value = 1.042022527221976
$ws.Range("J3").Value = 1.051590292994987
$ws.Range("K3").Value = 1.050627933206289
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.048929136109861
$ws.Range("N3").Value = 1.053083671830132

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047764708138949
$ws.Range("D4").Value = 1.04873868275095
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.047634105754515
$ws.Range("I4").Value = 1.042294370644415
$ws.Range("J4").Value = 1.052239922375405
$ws.Range("K4").Value = 1.051190853804473
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.050089000700849
$ws.Range("N4").Value = 1.053734223758783

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048133864169676
$ws.Range("D5").Value = 1.049026503752941
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.048173480392591
$ws.Range("I5").Value = 1.042408103201223
$ws.Range("J5").Value = 1.052512393994282
$ws.Range("K5").Value = 1.051426898441041
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.050575931016366
$ws.Range("N5").Value = 1.054007082318613

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048195816562565
$ws.Range("D6").Value = 1.049074803061212
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.048264010087945
$ws.Range("I6").Value = 1.042427167212126
$ws.Range("J6").Value = 1.052558106278151
$ws.Range("K6").Value = 1.051466495931888
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.050657649517191
$ws.Range("N6").Value = 1.054052859519169

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047769642871539
$ws.Range("D7").Value = 1.048742530449616
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.047641315173543
$ws.Range("I7").Value = 1.042295892503312
$ws.Range("J7").Value = 1.052243565631826
$ws.Range("K7").Value = 1.05119401022334
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.050095509727894
$ws.Range("N7").Value = 1.053737872189047

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045983829999774
$ws.Range("D8").Value = 1.047349722124808
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.045033562046783
$ws.Range("I8").Value = 1.041742505968483
$ws.Range("J8").Value = 1.050923450603902
$ws.Range("K8").Value = 1.050049896335301
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.047740072672029
$ws.Range("N8").Value = 1.052415882446364

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042819557270246
$ws.Range("D9").Value = 1.044880081492271
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.040418427378897
$ws.Range("I9").Value = 1.040749516501381
$ws.Range("J9").Value = 1.048576483126335
$ws.Range("K9").Value = 1.048013960838179
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.043566658727125
$ws.Range("N9").Value = 1.050065582005778

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040697776798541
$ws.Range("D10").Value = 1.043222968567604
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.037327197545815
$ws.Range("I10").Value = 1.04007532011024
$ws.Range("J10").Value = 1.046997454677648
$ws.Range("K10").Value = 1.046642957282836
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.040768060892111
$ws.Range("N10").Value = 1.048484311155577

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039775992458743
$ws.Range("D11").Value = 1.042502801128097
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.035984989676868
$ws.Range("I11").Value = 1.039780445825944
$ws.Range("J11").Value = 1.046310210084931
$ws.Range("K11").Value = 1.046045963945716
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.039552148878867
$ws.Range("N11").Value = 1.04779609059671

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039433133559911
$ws.Range("D12").Value = 1.042234896847177
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.035485859676985
$ws.Range("I12").Value = 1.039670470284591
$ws.Range("J12").Value = 1.046054400372432
$ws.Range("K12").Value = 1.045823705188535
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.039099870159579
$ws.Range("N12").Value = 1.047539917605078

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039506699267289
$ws.Range("D13").Value = 1.042292381537395
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.035592951074982
$ws.Range("I13").Value = 1.039694080663026
$ws.Range("J13").Value = 1.046109296836516
$ws.Range("K13").Value = 1.045871403587703
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.039196914512514
$ws.Range("N13").Value = 1.047594892028437

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039747661226231
$ws.Range("D14").Value = 1.042480664321101
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.035943743269652
$ws.Range("I14").Value = 1.039771364336944
$ws.Range("J14").Value = 1.046289075767101
$ws.Range("K14").Value = 1.046027602389225
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.039514776439135
$ws.Range("N14").Value = 1.047774926265725

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039896063623031
$ws.Range("D15").Value = 1.042596618062403
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.036159801206095
$ws.Range("I15").Value = 1.03981892212629
$ws.Range("J15").Value = 1.046399772159005
$ws.Range("K15").Value = 1.046123773984868
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.039710536972634
$ws.Range("N15").Value = 1.047885779859201

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040758887657631
$ws.Range("D16").Value = 1.043270707738627
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.037416196069566
$ws.Range("I16").Value = 1.040094827605239
$ws.Range("J16").Value = 1.047042990068766
$ws.Range("K16").Value = 1.046682506796385
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.0408486689731
$ws.Range("N16").Value = 1.048529911212173

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041299293256832
$ws.Range("D17").Value = 1.04369283796149
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.038203298325264
$ws.Range("I17").Value = 1.040267105156239
$ws.Range("J17").Value = 1.047445516740084
$ws.Range("K17").Value = 1.047032086000965
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.041561478045885
$ws.Range("N17").Value = 1.048933009517535

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041614210159196
$ws.Range("D18").Value = 1.043938806454663
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.038662047613236
$ws.Range("I18").Value = 1.040367308071248
$ws.Range("J18").Value = 1.04767996517772
$ws.Range("K18").Value = 1.047235667882721
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.04197685314038
$ws.Range("N18").Value = 1.049167790898844

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041721539381037
$ws.Range("D19").Value = 1.044022632676784
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.038818409817154
$ws.Range("I19").Value = 1.040401426680771
$ws.Range("J19").Value = 1.047759848863336
$ws.Range("K19").Value = 1.04730502968144
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.042118418985083
$ws.Range("N19").Value = 1.049247788028456

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041241343185282
$ws.Range("D20").Value = 1.043647573598273
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.038118886481553
$ws.Range("I20").Value = 1.040248650775101
$ws.Range("J20").Value = 1.047402364517272
$ws.Range("K20").Value = 1.046994612792471
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.041485041306244
$ws.Range("N20").Value = 1.048889796013617

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039676716915874
$ws.Range("D21").Value = 1.042425230892024
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.035840459729954
$ws.Range("I21").Value = 1.039748618568041
$ws.Range("J21").Value = 1.04623615023011
$ws.Range("K21").Value = 1.045981619839452
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.03942119170128
$ws.Range("N21").Value = 1.047721925568401

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038690266187285
$ws.Range("D22").Value = 1.04165436577123
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.034404589677703
$ws.Range("I22").Value = 1.039431645454328
$ws.Range("J22").Value = 1.045499796932231
$ws.Range("K22").Value = 1.045341763339521
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.038119884992055
$ws.Range("N22").Value = 1.046984526564377

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039213462619263
$ws.Range("D23").Value = 1.042063239286394
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.035166094403414
$ws.Range("I23").Value = 1.039599925080307
$ws.Range("J23").Value = 1.04589044918861
$ws.Range("K23").Value = 1.045681245165664
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.038810087685142
$ws.Range("N23").Value = 1.047375733591769

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041267529232196
$ws.Range("D24").Value = 1.043668027397281
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.038157029654218
$ws.Range("I24").Value = 1.040256990392129
$ws.Range("J24").Value = 1.04742186419548
$ws.Range("K24").Value = 1.047011546316648
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.041519580997039
$ws.Range("N24").Value = 1.048909323383604

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043639719675618
$ws.Range("D25").Value = 1.045520398424685
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.041614022510185
$ws.Range("I25").Value = 1.041008364143329
$ws.Range("J25").Value = 1.049185733980537
$ws.Range("K25").Value = 1.048542689370429
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.044648391380145
$ws.Range("N25").Value = 1.050675698066075

